# Add the 2022 column (S) to the "domestic taxes" table:
#  - S4 gets the year 2022, formatted like the other year headers (R4)
#  - S5 gets the 76.1% value, formatted like the other data cells (R5)
#  - leave the cursor where the author left it after entering the data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year header cell: copy R4's formatting onto S4, then set the new year.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Data cell: copy R5's formatting (incl. the 0.0 number format) onto S5.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 76.1

# Clear marching-ants clipboard state and move the selection cursor.
$excel.CutCopyMode = $false
$ws.Range("P8").Select() | Out-Null
